$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmark")

# Insert a new column before column F (shifts F:K -> G:L)
$ws.Columns.Item(6).Insert()

# Set the header for the newly inserted column
$ws.Cells.Item(1, 6).Value = "storage_bucket_name"
